$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (G=5503)
$ws.Range("H5").Value = 175.75
$ws.Range("I5").Value = 175.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 175.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -60.75
$ws.Range("N5").ClearContents()

# Row 99 (G=19883)
$ws.Range("H99").Value = 375.33334
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

# Row 112 (G=27960)
$ws.Range("H112").Value = 1655.9474
$ws.Range("J112").Value = 1814
$ws.Range("L112").Value = 5442
$ws.Range("N112").Value = -7658

# Row 116 (G=27778)
$ws.Range("H116").Value = 8188.4346
$ws.Range("I116").Value = 8244.9
$ws.Range("J116").Value = 8145
$ws.Range("K116").Value = 8244.9
$ws.Range("L116").Value = 8145
$ws.Range("M116").Value = -4802.9
$ws.Range("N116").Value = -15029

# Row 125 (G=36228)
$ws.Range("H125").Value = 1449.75
$ws.Range("I125").Value = 1400
$ws.Range("J125").Value = 1466.3334
$ws.Range("K125").Value = 12600
$ws.Range("L125").Value = 13197.0006
$ws.Range("M125").Value = -10140
$ws.Range("N125").Value = -18117.0006

# Row 130 (G=34691)
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040

# Row 132 (G=44049)
$ws.Range("H132").Value = 12312.3125
$ws.Range("I132").Value = 1591.1818
$ws.Range("K132").Value = 4773.5454
$ws.Range("M132").Value = -2243.5454

$ws = $wb.Worksheets.Item("ARM")
# Row 21 (G=3045)
$ws.Range("H21").Value = 6476.1665
$ws.Range("I21").Value = 5801.875
$ws.Range("J21").Value = 7824.75
$ws.Range("K21").Value = 5801.875
$ws.Range("L21").Value = 7824.75
$ws.Range("M21").Value = -5427.875
$ws.Range("N21").Value = -8572.75

# Row 32 (G=44147)
$ws.Range("H32").Value = 12663240
$ws.Range("I32").Value = 13703940
$ws.Range("K32").Value = 13703940
$ws.Range("M32").Value = -13703653

# Row 61 (G=43999)
$ws.Range("H61").Value = 7649.6665
$ws.Range("I61").Value = 7649.6665
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7649.6665
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -7437.6665
$ws.Range("N61").ClearContents()

# Row 101 (G=18518)
$ws.Range("H101").Value = 43666.668
$ws.Range("J101").Value = 43666.668
$ws.Range("L101").Value = 43666.668
$ws.Range("N101").Value = -50156.668

# Row 102 (G=19945)
$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 10000
$ws.Range("K102").Value = 10000
$ws.Range("M102").Value = -8378

# Row 104 (G=18672)
$ws.Range("H104").Value = 68997.5
$ws.Range("J104").Value = 68997.5
$ws.Range("L104").Value = 68997.5
$ws.Range("N104").Value = -75985.5

# Row 131 (G=34706)
$ws.Range("H131").Value = 51211.93
$ws.Range("J131").Value = 51211.93
$ws.Range("L131").Value = 51211.93
$ws.Range("N131").Value = -61291.93

# Row 136 (G=43999)
$ws.Range("H136").Value = 7649.6665
$ws.Range("I136").Value = 7649.6665
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 22948.9995
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -20398.9995
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (G=19947)
$ws.Range("H105").Value = 7252
$ws.Range("I105").Value = 6905.3335
$ws.Range("K105").Value = 6905.3335
$ws.Range("M105").Value = -5158.3335

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (G=44023)
$ws.Range("H31").Value = 6391.304
$ws.Range("I31").Value = 11870.333
$ws.Range("J31").Value = 2869.0715
$ws.Range("K31").Value = 11870.333
$ws.Range("L31").Value = 2869.0715
$ws.Range("M31").Value = -11575.333
$ws.Range("N31").Value = -3459.0715

# Row 34 (G=44023)
$ws.Range("H34").Value = 6391.304
$ws.Range("I34").Value = 11870.333
$ws.Range("J34").Value = 2869.0715
$ws.Range("K34").Value = 11870.333
$ws.Range("L34").Value = 2869.0715
$ws.Range("M34").Value = -11668.333
$ws.Range("N34").Value = -3273.0715

# Row 99 (G=36198)
$ws.Range("H99").Value = 28056780
$ws.Range("I99").Value = 6103311
$ws.Range("K99").Value = 6103311
$ws.Range("M99").Value = -6101813

# Row 126 (G=36198)
$ws.Range("H126").Value = 28056780
$ws.Range("I126").Value = 6103311
$ws.Range("K126").Value = 18309933
$ws.Range("M126").Value = -18307463

$ws = $wb.Worksheets.Item("CUL")
# Row 6 (G=4639)
$ws.Range("H6").Value = 416.45456
$ws.Range("I6").Value = 328.1
$ws.Range("J6").Value = 1300
$ws.Range("K6").Value = 984.3000000000001
$ws.Range("L6").Value = 3900
$ws.Range("M6").Value = -871.3000000000001
$ws.Range("N6").Value = -4126

# Row 54 (G=4671)
$ws.Range("H54").Value = 30674.908
$ws.Range("I54").Value = 2500
$ws.Range("K54").Value = 7500
$ws.Range("M54").Value = -6941

# Row 116 (G=27866)
$ws.Range("H116").Value = 3124.2
$ws.Range("I116").Value = 874
$ws.Range("K116").Value = 2622
$ws.Range("M116").Value = 820

# Row 136 (G=44093)
$ws.Range("H136").Value = 2725
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 6 (G=2108)
$ws.Range("H6").Value = 787.8570999999999
$ws.Range("I6").Value = 801.75
$ws.Range("K6").Value = 801.75
$ws.Range("M6").Value = -688.75

# Row 16 (G=2108)
$ws.Range("H16").Value = 787.8570999999999
$ws.Range("I16").Value = 801.75
$ws.Range("K16").Value = 801.75
$ws.Range("M16").Value = -551.75

# Row 97 (G=19940)
$ws.Range("H97").Value = 10556.765
$ws.Range("I97").Value = 3628.75
$ws.Range("J97").Value = 27184
$ws.Range("K97").Value = 3628.75
$ws.Range("L97").Value = 27184
$ws.Range("M97").Value = -3132.75
$ws.Range("N97").Value = -28176

# Row 107 (G=27802)
$ws.Range("H107").Value = 210.6
$ws.Range("I107").Value = 198.66667
$ws.Range("K107").Value = 198.66667
$ws.Range("M107").Value = 1721.33333

# Row 122 (G=36182)
$ws.Range("H122").Value = 5658.6
$ws.Range("I122").Value = 5067.5
$ws.Range("J122").Value = 6249.7
$ws.Range("K122").Value = 15202.5
$ws.Range("L122").Value = 18749.1
$ws.Range("M122").Value = -12752.5
$ws.Range("N122").Value = -23649.1

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (G=5277)
$ws.Range("H22").Value = 2918.9565
$ws.Range("J22").Value = 3175.5
$ws.Range("L22").Value = 3175.5
$ws.Range("N22").Value = -3765.5

# Row 27 (G=5277)
$ws.Range("H27").Value = 2918.9565
$ws.Range("J27").Value = 3175.5
$ws.Range("L27").Value = 3175.5
$ws.Range("N27").Value = -3389.5

# Row 46 (G=5282)
$ws.Range("H46").Value = 4349.9062
$ws.Range("I46").Value = 650
$ws.Range("K46").Value = 650
$ws.Range("M46").Value = -462

# Row 125 (G=34271)
$ws.Range("H125").Value = 50914.332
$ws.Range("J125").Value = 50914.332
$ws.Range("L125").Value = 50914.332
$ws.Range("N125").Value = -60754.332

# Row 140 (G=42503)
$ws.Range("H140").Value = 75994
$ws.Range("J140").Value = 75994
$ws.Range("L140").Value = 75994
$ws.Range("N140").Value = -86354

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (G=44029)
$ws.Range("H132").Value = 5123.0293
$ws.Range("I132").Value = 4579.913
$ws.Range("J132").Value = 6258.636
$ws.Range("K132").Value = 13739.739
$ws.Range("L132").Value = 18775.908
$ws.Range("M132").Value = -11209.739
$ws.Range("N132").Value = -23835.908

# Row 138 (G=42347)
$ws.Range("H138").Value = 75899.8
$ws.Range("J138").Value = 75899.8
$ws.Range("L138").Value = 75899.8
$ws.Range("N138").Value = -86179.8
